$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") values recomputed - regenerated save_data using K (strikeouts)
# instead of Strike# for rows 2-21.
$kValues = @{
    2  = 3
    3  = 0
    4  = 4
    5  = 1
    6  = 7
    7  = 2
    8  = 7
    9  = 0
    10 = 2
    11 = 2
    12 = 4
    13 = 4
    14 = 5
    15 = 1
    16 = 4
    17 = 1
    18 = 3
    19 = 0
    20 = 4
    21 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
